$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 becomes blank (was 936.0092711364777)
$ws.Range("D3").Value = $null

# C4 updated value
$ws.Range("C4").Value = 18.69152608107289

# C5 updated value
$ws.Range("C5").Value = 0

# Row 7 "Other" renamed to "Biogas" and D7 gets a new value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 411.8440793000502

# New row 8: "Other" with D8 = 0, matching style/format of A7 (bold/border/alignment)
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 0
